$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1  = 11.1
    2  = 8.6
    3  = 3.3
    4  = 10
    5  = 22.5
    6  = 13
    7  = 3.1
    8  = 9.699999999999999
    9  = 25.1
    11 = 16.5
    12 = 5.2
    13 = 12.5
    14 = 0.8
    15 = 2.8
    16 = 12
}

foreach ($row in $values.Keys) {
    $ws.Range("AO$row").Value = $values[$row]
}
